$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 8: Q3 2021 -> Q4 2021 report data refresh.
# Clear the old "Nota" text first so the stale shared string is dropped
# before we add the new strings (keeps shared-string ordering sane).
$ws.Range("M8").Value = ""

$ws.Range("A8").Value = 2021
$ws.Range("B8").Value = 44470
$ws.Range("C8").Value = 44561

$ws.Range("E8").Value = "Dirección de Planeación"
$ws.Range("F8").Value = "Decreto de Creación de la Universidad Politécnica de Pachuca en su artículo 26 de las facultades y obligaciones del Rector fracción X."
$ws.Range("G8").Value = "Trimestral"
$ws.Range("D8").Value = "Tercera Sesión Ordinaria 2021"

$ws.Range("H8").Value = 44530

$ws.Hyperlinks.Add($ws.Range("I8"), "https://drive.google.com/file/d/1zxT-oMD3k_jmwqtT2mCP6dQhO7d53Ar8/view?usp=sharing")

$ws.Range("K8").Value = 44571
$ws.Range("L8").Value = 44571

# Column M width shrinks now that the long "Nota" text is gone.
$ws.Columns("M").ColumnWidth = 22.7109375

# Scroll/selection bookkeeping recorded by Excel when the sheet was last saved.
$ws.Application.ActiveWindow.ScrollColumn = 6
$ws.Range("H12").Select()
